# Applies the "add DailyPriceData sheet / clear StockStreak" edit.

$wb = $excel.ActiveWorkbook

# The original "StockStreak" sheet holds the notes content that needs to
# move onto a brand-new second sheet named "DailyPriceData".
$src = $wb.Worksheets.Item("StockStreak")

# Capture the existing notes (A1:B3) before anything else changes.
$values = $src.Range("A1:B3").Value2

# Create the new sheet immediately after StockStreak.
$new = $wb.Worksheets.Add($null, $src)
$new.Name = "DailyPriceData"

# Re-fetch both sheets by name -- worksheet variables can go stale across
# structural operations like Add/Move.
$src = $wb.Worksheets.Item("StockStreak")
$new = $wb.Worksheets.Item("DailyPriceData")

# Write the original notes onto the new sheet.
$new.Range("A1:B3").Value = $values

# Add the new "Volume Score" row.
$new.Range("A4").Value = "Volume Score"
$new.Range("B4").Value = "100 / Average Daily Volume * Volume - 100"

# Column widths on the new sheet.
$new.Columns.Item(1).ColumnWidth = 11.5
$new.Columns.Item(2).ColumnWidth = 40.19921875

# Clear the old sheet's content -- it becomes empty.
$src.Range("A1:B3").Clear()

# Selection / active-cell bookkeeping to match the authored file.
$src.Range("A1:B3").Select()
$new.Range("A5").Select()

# Make the new sheet the active tab.
$new.Activate()
